# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp title (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 15:03"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 3480089
$ws.Range("C4").Value = 606
$ws.Range("E4").Value = 1791492

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 1888889
$ws.Range("C5").Value = 930
$ws.Range("E5").Value = 602427
$ws.Range("G5").Value = 29
$ws.Range("H5").Value = 72950

# --- Row 25: Argentina ---
$ws.Range("D25").Value = 45467
$ws.Range("E25").Value = 55872
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = 1926

# --- Row 28: Irak ---
$ws.Range("B28").Value = 81757
$ws.Range("C28").Value = 2022
$ws.Range("D28").Value = 50782
$ws.Range("E28").Value = 27630
$ws.Range("G28").Value = 95
$ws.Range("H28").Value = 3345

# --- Row 40: Paises Bajos ---
$ws.Range("B40").Value = 51146
$ws.Range("C40").Value = 53
$ws.Range("H40").Value = 6135

# --- Row 70: Dinamarca ---
$ws.Range("B70").Value = 13061
$ws.Range("C70").Value = 24
$ws.Range("D70").Value = 12160
$ws.Range("E70").Value = 291

# --- Row 80: Senegal ---
$ws.Range("B80").Value = 8243
$ws.Range("C80").Value = 45
$ws.Range("D80").Value = 5580
$ws.Range("E80").Value = 2513

# --- Rows 100/101: swap Grecia/Croacia order, Croacia now first with updated data ---
$ws.Range("A100").Value = "Croacia"
$ws.Range("B100").Value = 3827
$ws.Range("C100").Value = 52
$ws.Range("D100").Value = 2558
$ws.Range("E100").Value = 1149
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 120

$ws.Range("A101").Value = "Grecia"
$ws.Range("B101").Value = 3826
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 1374
$ws.Range("E101").Value = 2259
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 193

# --- Row 140: Burkina Faso ---
$ws.Range("B140").Value = 1037
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 873
$ws.Range("E140").Value = 111

# --- Row 147: Namibia ---
$ws.Range("B147").Value = 864
$ws.Range("C147").Value = 3
$ws.Range("D147").Value = 29
$ws.Range("E147").Value = 833
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 2
